# Fixed an issue with drive team data list being null.
# Adds a new "Drive Team Data" worksheet and refreshes the weighted-average
# columns (E:G) on the "Per Member Data" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the weighted-average columns (E,F,G) on "Per Member Data"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Per Member Data")

$ws3.Range("E2").Value = 70.55813953490454
$ws3.Range("F2").Value = 34.57613168743349
$ws3.Range("G2").Value = 22.76131687240414

$ws3.Range("E3").Value = 81.75294117656986
$ws3.Range("F3").Value = 46.90140845063206
$ws3.Range("G3").Value = 33.42253521145376

$ws3.Range("E4").Value = 68.94666666681476

$ws3.Range("E5").Value = 60.564705882457396
$ws3.Range("F5").Value = 37.815602836996234
$ws3.Range("G5").Value = 22.578014184354227

$ws3.Range("E6").Value = 93.24090909078771
$ws3.Range("F6").Value = 48.051282051094276
$ws3.Range("G6").Value = 41.153846153805915

$ws3.Range("E7").Value = 71.02721088439371
$ws3.Range("F7").Value = 26.35955056179175
$ws3.Range("G7").Value = 28.71910112358349

$ws3.Range("F8").Value = 38.16666666666667

$ws3.Range("E9").Value = 89.20512820490019
$ws3.Range("F9").Value = 48.051282051094276
$ws3.Range("G9").Value = 41.153846153805915

$ws3.Range("E10").Value = 73.98329853883885
$ws3.Range("F10").Value = 32.9638009050467
$ws3.Range("G10").Value = 26.687782805497104

# ---------------------------------------------------------------------
# 1b. Update the saved selection on "Match Data" (it's no longer the
#     active/tabbed sheet once "Drive Team Data" is added below).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Match Data")
$ws1.Range("H30").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Drive Team Data" worksheet at the end of the workbook
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Drive Team Data"

$ws4.Range("A2").Value = 64.85714285714286
$ws4.Range("B2").Value = 37.08666666666666
$ws4.Range("C2").Value = 22.06
$ws4.Range("D2").Value = 0.0
$ws4.Range("E2").Value = 0.0
$ws4.Range("F2").Value = 0.0
$ws4.Range("G2").Value = 0.0
$ws4.Range("H2").Value = 0.0
$ws4.Range("K2").Value = "B+M"

$ws4.Range("A3").Value = 87.63333333333334
$ws4.Range("B3").Value = 47.040000000000006
$ws4.Range("C3").Value = 38.28
$ws4.Range("D3").Value = 0.0
$ws4.Range("E3").Value = 0.0
$ws4.Range("F3").Value = 0.0
$ws4.Range("G3").Value = 0.0
$ws4.Range("H3").Value = 0.0
$ws4.Range("K3").Value = "E+Z"

$ws4.Range("A4").Value = 71.3075
$ws4.Range("B4").Value = 27.790476190476188
$ws4.Range("C4").Value = 28.46666666666667
$ws4.Range("D4").Value = 0.0
$ws4.Range("E4").Value = 0.0
$ws4.Range("F4").Value = 0.0
$ws4.Range("G4").Value = 0.0
$ws4.Range("H4").Value = 0.0
$ws4.Range("K4").Value = "L+C"

# Match the saved selection on the new sheet.
$ws4.Range("K5").Select()
